# "Updated data in readme file." - add a new "Use of force definition" source
# row to the Sources sheet, with its Wikipedia link, mirroring the existing
# rows/hyperlinks (Law enforcement / Geographic coordinate system).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")

# New row 4: description in column A, URL (as a hyperlink) in column B.
$ws.Range("A4").Value = "Use of force definition"
$ws.Range("B4").Value = "https://en.wikipedia.org/wiki/Use_of_force"

$ws.Hyperlinks.Add($ws.Range("B4"), "https://en.wikipedia.org/wiki/Use_of_force")

# Match the look of the other hyperlink cells (B2/B3) in the column.
$ws.Range("B4").Style = $ws.Range("B3").Style

# Leave the selection where Excel would land after typing into B4 and
# pressing Enter.
$ws.Range("B5").Select()
